# "Firm v0.3 - complemento a matrix"
# Rework the LCD splash-screen header (rows 2-4) on sheet "Hoja1" to read
# "Arduino ... DMX-512" / "Tester & Controller" / "Firm v0.0   Hard v0.0",
# highlight R8, and move the "About" label from row 69 up to row 68
# (P:U columns), matching the matrix/exit menu block above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row 2: "Arduino" (E2:K2) + "DMX-512" (M2:S2) -----------------------
$ws.Range("D2").ClearContents()
$ws.Cells.Item(2, 5).Value  = "A"   # E2
$ws.Cells.Item(2, 6).Value  = "r"   # F2
$ws.Cells.Item(2, 7).Value  = "d"   # G2
$ws.Cells.Item(2, 8).Value  = "u"   # H2
$ws.Cells.Item(2, 9).Value  = "i"   # I2
$ws.Cells.Item(2, 10).Value = "n"   # J2
$ws.Cells.Item(2, 11).Value = "o"   # K2
$ws.Cells.Item(2, 13).Value = "D"   # M2
$ws.Cells.Item(2, 14).Value = "M"   # N2
$ws.Cells.Item(2, 15).Value = "X"   # O2
$ws.Cells.Item(2, 16).Value = "-"   # P2
$ws.Cells.Item(2, 17).Value = "5"   # Q2
$ws.Cells.Item(2, 18).Value = "1"   # R2
$ws.Cells.Item(2, 19).Value = "2"   # S2
$ws.Range("T2").ClearContents()
$ws.Range("U2").ClearContents()

# --- Row 3: "Tester & Controller" (C3:H3, J3, L3:U3) --------------------
$ws.Cells.Item(3, 12).Value = "C"   # L3
$ws.Cells.Item(3, 13).Value = "o"   # M3
$ws.Cells.Item(3, 14).Value = "n"   # N3
$ws.Cells.Item(3, 15).Value = "t"   # O3
$ws.Cells.Item(3, 16).Value = "r"   # P3
$ws.Cells.Item(3, 17).Value = "o"   # Q3
$ws.Cells.Item(3, 18).Value = "l"   # R3
$ws.Cells.Item(3, 19).Value = "l"   # S3
$ws.Cells.Item(3, 20).Value = "e"   # T3
$ws.Cells.Item(3, 21).Value = "r"   # U3

# --- Row 4: "Firm v0.0" (B4:J4) + "Hard v0.0" (M4:U4) -------------------
$ws.Cells.Item(4, 2).Value  = "F"   # B4
$ws.Cells.Item(4, 3).Value  = "i"   # C4
$ws.Cells.Item(4, 4).Value  = "r"   # D4
$ws.Cells.Item(4, 5).Value  = "m"   # E4
$ws.Range("F4").ClearContents()
$ws.Cells.Item(4, 7).Value  = "v"   # G4
$ws.Cells.Item(4, 8).Value  = "0"   # H4
$ws.Cells.Item(4, 9).Value  = "."   # I4
$ws.Cells.Item(4, 10).Value = "0"   # J4
$ws.Range("K4").ClearContents()
$ws.Cells.Item(4, 13).Value = "H"   # M4
$ws.Cells.Item(4, 14).Value = "a"   # N4
$ws.Cells.Item(4, 15).Value = "r"   # O4
$ws.Cells.Item(4, 16).Value = "d"   # P4
$ws.Cells.Item(4, 18).Value = "v"   # R4
$ws.Cells.Item(4, 19).Value = "0"   # S4
$ws.Cells.Item(4, 20).Value = "."   # T4
$ws.Cells.Item(4, 21).Value = "0"   # U4

# --- Row 8: highlight R8 (matches the other orange separators) ----------
$ws.Range("R8").Interior.Color = $ws.Range("F8").Interior.Color

# --- Move the "About" label (P69:U69) up to (P68:U68) -------------------
$ws.Range("P69:U69").Copy($ws.Range("P68:U68"))
$ws.Range("P69:U69").ClearContents()

# --- Selection / view: land back on the header instead of row 43/X54 ----
[void]$ws.Range("E2").Select()
